# Add a new "JaaForkJoin-2" (Java ForkJoin) benchmark row to the results sheet.
# This inserts a new row above the current row 19 ("Java (streams)"), pushing
# all subsequent rows down by one, and populates the new row with the
# corresponding label and timing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Macbook Pro 2016")

# Insert a new blank row at row 19; existing rows 19-25 shift down to 20-26.
[void]$ws.Rows.Item(19).Insert()

# Populate the new row with the JaaForkJoin-2 results.
$ws.Cells.Item(19, 1).Value2 = "JaaForkJoin-2"   # A19 - Version/label
$ws.Cells.Item(19, 7).Value2 = 0.83              # G19
$ws.Cells.Item(19, 8).Value2 = 0.41              # H19
$ws.Cells.Item(19, 9).Value2 = 2.836             # I19
$ws.Cells.Item(19, 10).Value2 = 19.347           # J19
$ws.Cells.Item(19, 11).Value2 = 138              # K19
$ws.Cells.Item(19, 12).Value2 = 1042             # L19

# Match the saved selection state from the edit (active cell I19).
[void]$ws.Range("I19").Select()
